$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText, $matchWholeWord) {
    $rng = $d.Content
    $result = $rng.Find.Execute($findText, $true, $matchWholeWord, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $findText"
    }
}

# 1. "neccesary" -> "necessary" (typo fix; MatchWholeWord so "uneccesary" below isn't touched)
Replace-Text "neccesary" "necessary" $true

# 2. "accesse" -> "access" (typo fix)
Replace-Text "accesse player information" "access player information" $false

# 3. "Github" -> "GitHub", "changeswhat" -> "changes what"
Replace-Text "Changes will be made through Github with explanations of the changeswhat " "Changes will be made through GitHub with explanations of the changes what " $false

# 4. "SPMP : " -> "SPMP: " (remove space before colon)
Replace-Text "SPMP : Software Project Management Plan" "SPMP: Software Project Management Plan" $false

# 5. "API : " -> "API: " (remove space before colon). The word "API" itself lives in the
#    same run as the preceding <w:tab/>; editing that run flattens the tab to a literal
#    character, so instead we only touch the following run (" : Application...") and
#    shift the colon there, leaving "API" (and its <w:tab/>) run untouched.
Replace-Text " : Application Program interface" ": Application Program interface" $false

# 6. "I will be in charge of determining" -> "I will oversee determining"
Replace-Text "I will be in charge of determining objectives for each prototype" "I will oversee determining objectives for each prototype" $false

# 7. "no client except myself" -> "no client except me"
Replace-Text "no client except myself" "no client except me" $false

# 8. "uneccesary" -> "unnecessary" (typo fix)
Replace-Text "uneccesary" "unnecessary" $true

# 9. "Time :Have to finish by December" -> "Time: Must finish by December"
Replace-Text "Time :Have to finish by December" "Time: Must finish by December" $false

# 10. "Need to take into account every" -> "Need to consider every"
Replace-Text "Need to take into account every new development idea for our program." "Need to consider every new development idea for our program." $false

# 11. Append "Firebase(If Needed)" right after "Database: "
$rng = $d.Content
$found = $rng.Find.Execute("Database: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.InsertAfter("Firebase(If Needed)")
} else {
    Write-Output "NOT FOUND: Database: "
}

# 12. " GOdaddy, Google Cloud, Heroku" -> " GoDaddy"
Replace-Text " GOdaddy, Google Cloud, Heroku" " GoDaddy" $false

# 13. "Acceptance(final test)" -> "Acceptance (final test)"
Replace-Text "Acceptance(final test)" "Acceptance (final test)" $false

Write-Output "done"
